$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5

$ws.Range("I1").Font.Color = 0
$ws.Range("K1").Font.Color = 0
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0
